# Auto-generated edit script: update Asura_Profits market-price derived columns (H-N)
# across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets, per scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 243.71428
$ws.Range("I33").Value = 243.71428
$ws.Range("K33").Value = 243.71428
$ws.Range("M33").Value = -14.71428
$ws.Range("H43").Value = 2000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 2000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 2000
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -2138
$ws.Range("H64").Value = 4279.3335
$ws.Range("I64").Value = 3599.4
$ws.Range("J64").Value = 4379.3237
$ws.Range("K64").Value = 3599.4
$ws.Range("L64").Value = 4379.3237
$ws.Range("M64").Value = -3351.4
$ws.Range("N64").Value = -4875.3237
$ws.Range("H67").Value = 4279.3335
$ws.Range("I67").Value = 3599.4
$ws.Range("J67").Value = 4379.3237
$ws.Range("K67").Value = 3599.4
$ws.Range("L67").Value = 4379.3237
$ws.Range("M67").Value = -2741.4
$ws.Range("N67").Value = -6095.3237
$ws.Range("H129").Value = 959.5472
$ws.Range("I129").Value = 419.30768
$ws.Range("J129").Value = 1135.125
$ws.Range("K129").Value = 1257.92304
$ws.Range("L129").Value = 3405.375
$ws.Range("M129").Value = 3742.07696
$ws.Range("N129").Value = -13405.375
$ws.Range("H131").Value = 4533.5835
$ws.Range("I131").Value = 866.73334
$ws.Range("J131").Value = 10645
$ws.Range("K131").Value = 2600.20002
$ws.Range("L131").Value = 31935
$ws.Range("M131").Value = 2439.79998
$ws.Range("N131").Value = -42015
$ws.Range("H138").Value = 3129802.8
$ws.Range("I138").Value = 7695083.5
$ws.Range("J138").Value = 6189.7104
$ws.Range("K138").Value = 23085250.5
$ws.Range("L138").Value = 18569.1312
$ws.Range("M138").Value = -23080110.5
$ws.Range("N138").Value = -28849.1312
$ws.Range("H141").Value = 18705.75
$ws.Range("I141").Value = 13255.571
$ws.Range("K141").Value = 39766.713
$ws.Range("M141").Value = -34586.713

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15288.72
$ws.Range("I32").Value = 13564.962
$ws.Range("K32").Value = 13564.962
$ws.Range("M32").Value = -13277.962
$ws.Range("H110").Value = 1702.5
$ws.Range("I110").Value = 1686.6666
$ws.Range("J110").Value = 1750
$ws.Range("K110").Value = 1686.6666
$ws.Range("L110").Value = 1750
$ws.Range("M110").Value = 358.3334
$ws.Range("N110").Value = -5840
$ws.Range("H133").Value = 63392
$ws.Range("J133").Value = 63392
$ws.Range("L133").Value = 63392
$ws.Range("N133").Value = -68452
$ws.Range("H135").Value = 44357
$ws.Range("J135").Value = 43851.5
$ws.Range("L135").Value = 43851.5
$ws.Range("N135").Value = -53991.5
$ws.Range("H138").Value = 62000
$ws.Range("J138").Value = 62000
$ws.Range("L138").Value = 62000
$ws.Range("N138").Value = -72280
$ws.Range("H140").Value = 46650.652
$ws.Range("J140").Value = 46650.652
$ws.Range("L140").Value = 46650.652
$ws.Range("N140").Value = -57010.652
$ws.Range("H141").Value = 56629.832
$ws.Range("J141").Value = 56629.832
$ws.Range("L141").Value = 56629.832
$ws.Range("N141").Value = -66989.83199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 74680
$ws.Range("J135").Value = 74680
$ws.Range("L135").Value = 74680
$ws.Range("N135").Value = -84820

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1344.2667
$ws.Range("I16").Value = 1108.125
$ws.Range("J16").Value = 1614.1428
$ws.Range("K16").Value = 1108.125
$ws.Range("L16").Value = 1614.1428
$ws.Range("M16").Value = -821.125
$ws.Range("N16").Value = -2188.1428
$ws.Range("H31").Value = 3170.558
$ws.Range("I31").Value = 1978.909
$ws.Range("K31").Value = 1978.909
$ws.Range("M31").Value = -1683.909
$ws.Range("H34").Value = 3170.558
$ws.Range("I34").Value = 1978.909
$ws.Range("K34").Value = 1978.909
$ws.Range("M34").Value = -1776.909
$ws.Range("H62").Value = 67307.5
$ws.Range("I62").Value = 102322
$ws.Range("K62").Value = 102322
$ws.Range("M62").Value = -101698
$ws.Range("H65").Value = 67307.5
$ws.Range("I65").Value = 102322
$ws.Range("K65").Value = 511610
$ws.Range("M65").Value = -508490
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H113").Value = 1344.2667
$ws.Range("I113").Value = 1108.125
$ws.Range("J113").Value = 1614.1428
$ws.Range("K113").Value = 1108.125
$ws.Range("L113").Value = 1614.1428
$ws.Range("M113").Value = 1061.875
$ws.Range("N113").Value = -5954.1428
$ws.Range("H134").Value = 1733.2609
$ws.Range("I134").Value = 1541.75
$ws.Range("J134").Value = 2171
$ws.Range("K134").Value = 4625.25
$ws.Range("L134").Value = 6513
$ws.Range("M134").Value = -2090.25
$ws.Range("N134").Value = -11583
$ws.Range("H137").Value = 56819
$ws.Range("J137").Value = 74698.336
$ws.Range("L137").Value = 74698.336
$ws.Range("N137").Value = -84898.336
$ws.Range("H138").Value = 58753.57
$ws.Range("J138").Value = 58753.57
$ws.Range("L138").Value = 58753.57
$ws.Range("N138").Value = -69033.57000000001
$ws.Range("H140").Value = 73180
$ws.Range("J140").Value = 73180
$ws.Range("L140").Value = 73180
$ws.Range("N140").Value = -83540

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 4626.593
$ws.Range("I113").Value = 566.6667
$ws.Range("J113").Value = 5134.0835
$ws.Range("K113").Value = 1700.0001
$ws.Range("L113").Value = 15402.2505
$ws.Range("M113").Value = 469.9999
$ws.Range("N113").Value = -19742.2505
$ws.Range("H114").Value = 476
$ws.Range("I114").Value = 318.44446
$ws.Range("J114").Value = 594.1667
$ws.Range("K114").Value = 955.33338
$ws.Range("L114").Value = 1782.5001
$ws.Range("M114").Value = 2298.66662
$ws.Range("N114").Value = -8290.500099999999
$ws.Range("H117").Value = 125207.25
$ws.Range("I117").Value = 276.33334
$ws.Range("K117").Value = 829.0000200000001
$ws.Range("M117").Value = 2612.99998
$ws.Range("H131").Value = 18185314
$ws.Range("I131").Value = 8796.154
$ws.Range("J131").Value = 23811380
$ws.Range("K131").Value = 26388.462
$ws.Range("L131").Value = 71434140
$ws.Range("M131").Value = -21348.462
$ws.Range("N131").Value = -71444220

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 10284
$ws.Range("J109").Value = 10284
$ws.Range("L109").Value = 10284
$ws.Range("N109").Value = -12364
$ws.Range("H126").Value = 3710
$ws.Range("I126").Value = 2800
$ws.Range("J126").Value = 3937.5
$ws.Range("K126").Value = 8400
$ws.Range("L126").Value = 11812.5
$ws.Range("M126").Value = -5930
$ws.Range("N126").Value = -16752.5
$ws.Range("H133").Value = 62310
$ws.Range("J133").Value = 62310
$ws.Range("L133").Value = 62310
$ws.Range("N133").Value = -72430
$ws.Range("H141").Value = 47710
$ws.Range("J141").Value = 47710
$ws.Range("L141").Value = 47710
$ws.Range("N141").Value = -58070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 11197.667
$ws.Range("I40").Value = 14444.75
$ws.Range("J40").Value = 8600
$ws.Range("K40").Value = 14444.75
$ws.Range("L40").Value = 8600
$ws.Range("M40").Value = -14308.75
$ws.Range("N40").Value = -8872
$ws.Range("H122").Value = 15529149
$ws.Range("I122").Value = 18633778
$ws.Range("J122").Value = 6001.6665
$ws.Range("K122").Value = 55901334
$ws.Range("L122").Value = 18004.9995
$ws.Range("M122").Value = -55898884
$ws.Range("N122").Value = -22904.9995
$ws.Range("H134").Value = 76660
$ws.Range("J134").Value = 76660
$ws.Range("L134").Value = 76660
$ws.Range("N134").Value = -86800
$ws.Range("H139").Value = 50175
$ws.Range("J139").Value = 50175
$ws.Range("L139").Value = 50175
$ws.Range("N139").Value = -60455
$ws.Range("H140").Value = 76187.5
$ws.Range("J140").Value = 76187.5
$ws.Range("L140").Value = 76187.5
$ws.Range("N140").Value = -86547.5
$ws.Range("H141").Value = 47000
$ws.Range("J141").Value = 47000
$ws.Range("L141").Value = 47000
$ws.Range("N141").Value = -57360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 20835512
$ws.Range("I122").Value = 25002014
$ws.Range("K122").Value = 75006042
$ws.Range("M122").Value = -75003592
$ws.Range("H132").Value = 2195.5186
$ws.Range("I132").Value = 1699.25
$ws.Range("K132").Value = 5097.75
$ws.Range("M132").Value = -2567.75
$ws.Range("H133").Value = 46552
$ws.Range("J133").Value = 46552
$ws.Range("L133").Value = 46552
$ws.Range("N133").Value = -56672
$ws.Range("H135").Value = 50313.125
$ws.Range("J135").Value = 50313.125
$ws.Range("L135").Value = 50313.125
$ws.Range("N135").Value = -60453.125
$ws.Range("H138").Value = 58474.75
$ws.Range("J138").Value = 58474.75
$ws.Range("L138").Value = 58474.75
$ws.Range("N138").Value = -68754.75
$ws.Range("H139").Value = 60787.145
$ws.Range("J139").Value = 60787.145
$ws.Range("L139").Value = 60787.145
$ws.Range("N139").Value = -71067.14499999999
$ws.Range("H140").Value = 42189.89
$ws.Range("J140").Value = 42189.89
$ws.Range("L140").Value = 42189.89
$ws.Range("N140").Value = -52549.89
$ws.Range("H141").Value = 89130.57000000001
$ws.Range("J141").Value = 89130.57000000001
$ws.Range("L141").Value = 89130.57000000001
$ws.Range("N141").Value = -99490.57000000001
